$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 283.1111
$ws.Range("I38").Value = 283.1111
$ws.Range("K38").Value = 849.3333
$ws.Range("M38").Value = -477.3333
$ws.Range("H55").Value = 1008.5455
$ws.Range("J55").Value = 2000
$ws.Range("L55").Value = 2000
$ws.Range("N55").Value = -2428
$ws.Range("H58").Value = 2658.75
$ws.Range("I58").Value = 80
$ws.Range("J58").Value = 5237.5
$ws.Range("K58").Value = 240
$ws.Range("L58").Value = 15712.5
$ws.Range("M58").Value = -90
$ws.Range("N58").Value = -16012.5
$ws.Range("H76").Value = 3625.75
$ws.Range("I76").Value = 3334.3333
$ws.Range("K76").Value = 3334.3333
$ws.Range("M76").Value = -3019.3333
$ws.Range("H79").Value = 3625.75
$ws.Range("I79").Value = 3334.3333
$ws.Range("K79").Value = 3334.3333
$ws.Range("M79").Value = -2242.3333
$ws.Range("H132").Value = 4744.1113
$ws.Range("I132").Value = 5054.45
$ws.Range("J132").Value = 3857.4285
$ws.Range("K132").Value = 15163.35
$ws.Range("L132").Value = 11572.2855
$ws.Range("M132").Value = -12633.35
$ws.Range("N132").Value = -16632.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 264
$ws.Range("I5").Value = 302.33334
$ws.Range("J5").Value = 149
$ws.Range("K5").Value = 302.33334
$ws.Range("L5").Value = 149
$ws.Range("M5").Value = -190.33334
$ws.Range("N5").Value = -373
$ws.Range("H88").Value = 1469.3529
$ws.Range("I88").Value = 1194.6666
$ws.Range("J88").Value = 1778.375
$ws.Range("K88").Value = 1194.6666
$ws.Range("L88").Value = 1778.375
$ws.Range("M88").Value = -788.6666
$ws.Range("N88").Value = -2590.375
$ws.Range("H91").Value = 1469.3529
$ws.Range("I91").Value = 1194.6666
$ws.Range("J91").Value = 1778.375
$ws.Range("K91").Value = 1194.6666
$ws.Range("L91").Value = 1778.375
$ws.Range("M91").Value = 209.3334
$ws.Range("N91").Value = -4586.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 264
$ws.Range("I4").Value = 302.33334
$ws.Range("J4").Value = 149
$ws.Range("K4").Value = 302.33334
$ws.Range("L4").Value = 149
$ws.Range("M4").Value = -187.33334
$ws.Range("N4").Value = -379
$ws.Range("H22").Value = 619.4286
$ws.Range("I22").Value = 1098.6666
$ws.Range("J22").Value = 260
$ws.Range("K22").Value = 1098.6666
$ws.Range("L22").Value = 260
$ws.Range("M22").Value = -925.6666
$ws.Range("N22").Value = -606
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H75").Value = 10000
$ws.Range("J75").Value = 10000
$ws.Range("L75").Value = 10000
$ws.Range("N75").Value = -11872
$ws.Range("H78").Value = 10000
$ws.Range("J78").Value = 10000
$ws.Range("L78").Value = 30000
$ws.Range("N78").Value = -39360
$ws.Range("H99").Value = 3090.5
$ws.Range("I99").Value = 3101
$ws.Range("K99").Value = 3101
$ws.Range("M99").Value = -1603

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 19157
$ws.Range("I47").Value = 19033
$ws.Range("J47").Value = 19250
$ws.Range("K47").Value = 19033
$ws.Range("L47").Value = 19250
$ws.Range("M47").Value = -18467
$ws.Range("N47").Value = -20382

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 63.869564
$ws.Range("I2").Value = 90.333336
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 542.000016
$ws.Range("L2").Value = 210
$ws.Range("M2").Value = -429.000016
$ws.Range("N2").Value = -436
$ws.Range("H26").Value = 600.7692
$ws.Range("I26").Value = 381.1
$ws.Range("J26").Value = 1333
$ws.Range("K26").Value = 1143.3
$ws.Range("L26").Value = 3999
$ws.Range("M26").Value = -855.3000000000002
$ws.Range("N26").Value = -4575
$ws.Range("H40").Value = 150
$ws.Range("J40").Value = 150
$ws.Range("L40").Value = 600
$ws.Range("N40").Value = -738
$ws.Range("H49").Value = 2278.4285
$ws.Range("J49").Value = 1990
$ws.Range("L49").Value = 5970
$ws.Range("N49").Value = -6282
$ws.Range("H110").Value = 900
$ws.Range("I110").Value = 900
$ws.Range("K110").Value = 2700
$ws.Range("M110").Value = 1390
$ws.Range("H122").Value = 1996.5
$ws.Range("I122").Value = 1995
$ws.Range("K122").Value = 17955
$ws.Range("M122").Value = -15505
$ws.Range("H140").Value = 669523.1
$ws.Range("I140").Value = 772295.9399999999
$ws.Range("K140").Value = 2316887.82
$ws.Range("M140").Value = -2311707.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7559.8
$ws.Range("I70").Value = 7099.7144
$ws.Range("K70").Value = 7099.7144
$ws.Range("M70").Value = -6829.7144
$ws.Range("H73").Value = 7559.8
$ws.Range("I73").Value = 7099.7144
$ws.Range("K73").Value = 7099.7144
$ws.Range("M73").Value = -6163.7144
$ws.Range("H100").Value = 25000
$ws.Range("J100").Value = 25000
$ws.Range("L100").Value = 25000
$ws.Range("N100").Value = -27164
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3273.875
$ws.Range("I40").Value = 2238.4
$ws.Range("J40").Value = 4999.6665
$ws.Range("K40").Value = 2238.4
$ws.Range("L40").Value = 4999.6665
$ws.Range("M40").Value = -2102.4
$ws.Range("N40").Value = -5271.6665
$ws.Range("H46").Value = 2866.3333
$ws.Range("I46").Value = 2175
$ws.Range("K46").Value = 2175
$ws.Range("M46").Value = -1987
$ws.Range("H98").Value = 54676.332
$ws.Range("J98").Value = 54676.332
$ws.Range("L98").Value = 54676.332
$ws.Range("N98").Value = -60666.332
$ws.Range("H100").Value = 1660.5
$ws.Range("I100").Value = 1752.6
$ws.Range("J100").Value = 1200
$ws.Range("K100").Value = 1752.6
$ws.Range("L100").Value = 1200
$ws.Range("M100").Value = -1211.6
$ws.Range("N100").Value = -2282
$ws.Range("H122").Value = 3587.8333
$ws.Range("I122").Value = 3666.6667
$ws.Range("J122").Value = 3509
$ws.Range("K122").Value = 11000.0001
$ws.Range("L122").Value = 10527
$ws.Range("M122").Value = -8550.000100000001
$ws.Range("N122").Value = -15427

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 68
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 68
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 68
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -298
$ws.Range("H122").Value = 4904.25
$ws.Range("I122").Value = 4865.8
$ws.Range("J122").Value = 4968.3335
$ws.Range("K122").Value = 14597.4
$ws.Range("L122").Value = 14905.0005
$ws.Range("M122").Value = -12147.4
$ws.Range("N122").Value = -19805.0005
